$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.38281187720569
$ws.Range("C2").Value = 6.832164312212983
$ws.Range("D2").Value = 7.896926053294598
$ws.Range("E2").Value = 12.68668211513331
$ws.Range("F2").Value = 37.45241620910736
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 28.21906590238764
$ws.Range("J2").Value = 10.03606963874178
$ws.Range("K2").Value = 10.72411960809373
$ws.Range("L2").Value = 11.06668315814925
$ws.Range("O2").Value = 28.97586259261839
# Row 3
$ws.Range("B3").Value = 13.15373421135726
$ws.Range("C3").Value = 6.775562307236552
$ws.Range("D3").Value = 7.876471785846572
$ws.Range("E3").Value = 12.703939935333
$ws.Range("F3").Value = 37.55689932773253
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 28.3246128628004
$ws.Range("J3").Value = 10.0558393029501
$ws.Range("K3").Value = 10.56022031431843
$ws.Range("L3").Value = 11.06398152686172
$ws.Range("O3").Value = 29.07451256594343
# Row 4
$ws.Range("B4").Value = 13.01292755549217
$ws.Range("C4").Value = 6.740177643211755
$ws.Range("D4").Value = 7.864896692918354
$ws.Range("E4").Value = 12.71604418821457
$ws.Range("F4").Value = 37.62830855442292
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 28.39395801630201
$ws.Range("J4").Value = 10.06870283130541
$ws.Range("K4").Value = 10.45963443465757
$ws.Range("L4").Value = 11.06363899743057
$ws.Range("O4").Value = 29.14017482780638
# Row 5
$ws.Range("B5").Value = 12.95558083655272
$ws.Range("C5").Value = 6.725605647484099
$ws.Range("D5").Value = 7.860430381992954
$ws.Range("E5").Value = 12.72135647915258
$ws.Range("F5").Value = 37.65923092040819
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 28.4233584200218
$ws.Range("J5").Value = 10.07412755763607
$ws.Range("K5").Value = 10.41870426955324
$ws.Range("L5").Value = 11.06383173903096
$ws.Range("O5").Value = 29.16821213268886
# Row 6
$ws.Range("B6").Value = 12.94606256384295
$ws.Range("C6").Value = 6.723176951550812
$ws.Range("D6").Value = 7.859703982556993
$ws.Range("E6").Value = 12.72226152972425
$ws.Range("F6").Value = 37.66447556738143
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 28.4283093007196
$ws.Range("J6").Value = 10.07503937950682
$ws.Range("K6").Value = 10.41191281736044
$ws.Range("L6").Value = 11.06388385230211
$ws.Range("O6").Value = 29.17294495724375
# Row 7
$ws.Range("B7").Value = 13.012153928391
$ws.Range("C7").Value = 6.739981729303392
$ws.Range("D7").Value = 7.864835439759278
$ws.Range("E7").Value = 12.71611429350102
$ws.Range("F7").Value = 37.62871820782303
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 28.39434989729133
$ws.Range("J7").Value = 10.06877525056701
$ws.Range("K7").Value = 10.45908213316095
$ws.Range("L7").Value = 11.06364024954631
$ws.Range("O7").Value = 29.14054776942824
# Row 8
$ws.Range("B8").Value = 13.3038999270391
$ws.Range("C8").Value = 6.812782541872599
$ws.Range("D8").Value = 7.88967119748086
$ws.Range("E8").Value = 12.69231998383495
$ws.Range("F8").Value = 37.48693475803132
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 28.25451665616574
$ws.Range("J8").Value = 10.04273604730381
$ws.Range("K8").Value = 10.66762508636207
$ws.Range("L8").Value = 11.06547926912064
$ws.Range("O8").Value = 29.00882008769197
# Row 9
$ws.Range("B9").Value = 13.87173661705192
$ws.Range("C9").Value = 6.950273706439337
$ws.Range("D9").Value = 7.946032235522976
$ws.Range("E9").Value = 12.65759939393346
$ws.Range("F9").Value = 37.26656081265286
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 28.01630222525489
$ws.Range("J9").Value = 9.997404019971057
$ws.Range("K9").Value = 11.07492515923822
$ws.Range("L9").Value = 11.07946813859805
$ws.Range("O9").Value = 28.79092703572812
# Row 10
$ws.Range("B10").Value = 14.28233424825462
$ws.Range("C10").Value = 7.047726702870719
$ws.Range("D10").Value = 7.991899922352068
$ws.Range("E10").Value = 12.63933534555403
$ws.Range("F10").Value = 37.13991799865087
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 27.86320918337089
$ws.Range("J10").Value = 9.967563316600048
$ws.Range("K10").Value = 11.37049019183829
$ws.Range("L10").Value = 11.09598834228598
$ws.Range("O10").Value = 28.655526989793
# Row 11
$ws.Range("B11").Value = 14.46689614775682
$ws.Range("C11").Value = 7.091214984104803
$ws.Range("D11").Value = 8.013686470472233
$ws.Range("E11").Value = 12.63259190499313
$ws.Range("F11").Value = 37.08998309681482
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 27.79831942979865
$ws.Range("J11").Value = 9.954734229075743
$ws.Range("K11").Value = 11.50360955901336
$ws.Range("L11").Value = 11.10483874285126
$ws.Range("O11").Value = 28.5992991658652
# Row 12
$ws.Range("B12").Value = 14.53640271314375
$ws.Range("C12").Value = 7.107555553865816
$ws.Range("D12").Value = 8.022064575736092
$ws.Range("E12").Value = 12.63026264678017
$ws.Range("F12").Value = 37.07217913818938
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 27.77443069409023
$ws.Range("J12").Value = 9.949982941436955
$ws.Range("K12").Value = 11.55378327987766
$ws.Range("L12").Value = 11.10838022199987
$ws.Range("O12").Value = 28.57877960485519
# Row 13
$ws.Range("B13").Value = 14.52145124534587
$ws.Range("C13").Value = 7.104042087601346
$ws.Range("D13").Value = 8.02025457979277
$ws.Range("E13").Value = 12.63075432845091
$ws.Range("F13").Value = 37.07596435793606
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 27.77954515492711
$ws.Range("J13").Value = 9.951001472298387
$ws.Range("K13").Value = 11.54298863155069
$ws.Range("L13").Value = 11.10760908248935
$ws.Range("O13").Value = 28.58316447759692
# Row 14
$ws.Range("B14").Value = 14.47262252435043
$ws.Range("C14").Value = 7.092561922901125
$ws.Range("D14").Value = 8.014373202648663
$ws.Range("E14").Value = 12.63239578366404
$ws.Range("F14").Value = 37.08849619428577
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 27.79634038784782
$ws.Range("J14").Value = 9.954341199477955
$ws.Range("K14").Value = 11.50774234916403
$ws.Range("L14").Value = 11.10512630585988
$ws.Range("O14").Value = 28.59759551749138
# Row 15
$ws.Range("B15").Value = 14.44266180715475
$ws.Range("C15").Value = 7.085513204794651
$ws.Range("D15").Value = 8.01078722592873
$ws.Range("E15").Value = 12.63343041598829
$ws.Range("F15").Value = 37.09631629193385
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 27.80671698276855
$ws.Range("J15").Value = 9.956400776542603
$ws.Range("K15").Value = 11.48612102088526
$ws.Range("L15").Value = 11.10363021948705
$ws.Range("O15").Value = 28.60653560382215
# Row 16
$ws.Range("B16").Value = 14.27022249013736
$ws.Range("C16").Value = 7.044867134276743
$ws.Range("D16").Value = 7.99049427899344
$ws.Range("E16").Value = 12.63980747717007
$ws.Range("F16").Value = 37.14333595302956
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 27.86754553467521
$ws.Range("J16").Value = 9.968416696118151
$ws.Range("K16").Value = 11.36175988882503
$ws.Range("L16").Value = 11.09543664962559
$ws.Range("O16").Value = 28.65930968274095
# Row 17
$ws.Range("B17").Value = 14.16382143514836
$ws.Range("C17").Value = 7.01971158130447
$ws.Range("D17").Value = 7.978278019207361
$ws.Range("E17").Value = 12.64411993560567
$ws.Range("F17").Value = 37.17414807369154
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 27.90607932809586
$ws.Range("J17").Value = 9.975978748425268
$ws.Range("K17").Value = 11.28509481925736
$ws.Range("L17").Value = 11.09075077724929
$ws.Range("O17").Value = 28.69306005390411
# Row 18
$ws.Range("B18").Value = 14.10241692470874
$ws.Range("C18").Value = 7.005163734280216
$ws.Range("D18").Value = 7.971338524420103
$ws.Range("E18").Value = 12.64674769446148
$ws.Range("F18").Value = 37.19259276978354
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 27.92869031617543
$ws.Range("J18").Value = 9.98039844417826
$ws.Range("K18").Value = 11.24087598551031
$ws.Range("L18").Value = 11.08818139574677
$ws.Range("O18").Value = 28.71297741995767
# Row 19
$ws.Range("B19").Value = 14.08159311839188
$ws.Range("C19").Value = 7.000224699055282
$ws.Range("D19").Value = 7.969004001555543
$ws.Range("E19").Value = 12.64766273617574
$ws.Range("F19").Value = 37.19896184474565
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 27.93642284101615
$ws.Range("J19").Value = 9.981906947161551
$ws.Range("K19").Value = 11.22588446629315
$ws.Range("L19").Value = 11.08733311378244
$ws.Range("O19").Value = 28.71980780832604
# Row 20
$ws.Range("B20").Value = 14.17516974204164
$ws.Range("C20").Value = 7.022397658286919
$ws.Range("D20").Value = 7.979569491059434
$ws.Range("E20").Value = 12.64364562217031
$ws.Range("F20").Value = 37.17079329519781
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 27.90193103813732
$ws.Range("J20").Value = 9.975166492231789
$ws.Range("K20").Value = 11.29326901135258
$ws.Range("L20").Value = 11.09123658972175
$ws.Range("O20").Value = 28.68941499149013
# Row 21
$ws.Range("B21").Value = 14.48697558996007
$ws.Range("C21").Value = 7.095937434564773
$ws.Range("D21").Value = 8.016097267604298
$ws.Range("E21").Value = 12.63190756587277
$ws.Range("F21").Value = 37.08478528113556
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 27.79138866682548
$ws.Range("J21").Value = 9.953357345093952
$ws.Range("K21").Value = 11.5181017732066
$ws.Range("L21").Value = 11.1058504166956
$ws.Range("O21").Value = 28.59333579249629
# Row 22
$ws.Range("B22").Value = 14.68849893577594
$ws.Range("C22").Value = 7.143253988765075
$ws.Range("D22").Value = 8.040714271569176
$ws.Range("E22").Value = 12.62554331516855
$ws.Range("F22").Value = 37.03501720018937
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 27.72312749127765
$ws.Range("J22").Value = 9.939726187723176
$ws.Range("K22").Value = 11.66365017771092
$ws.Range("L22").Value = 11.11650800842922
$ws.Range("O22").Value = 28.53504679938781
# Row 23
$ws.Range("B23").Value = 14.58116910890288
$ws.Range("C23").Value = 7.118070519978168
$ws.Range("D23").Value = 8.027509143517687
$ws.Range("E23").Value = 12.6288206657584
$ws.Range("F23").Value = 37.06098935579058
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 27.75919512094058
$ws.Range("J23").Value = 9.946944578342016
$ws.Range("K23").Value = 11.58610968233005
$ws.Range("L23").Value = 11.11071928027739
$ws.Range("O23").Value = 28.56574425533531
# Row 24
$ws.Range("B24").Value = 14.17003989375941
$ws.Range("C24").Value = 7.021183549003175
$ws.Range("D24").Value = 7.978985355298018
$ws.Range("E24").Value = 12.64385959676397
$ws.Range("F24").Value = 37.17230771553382
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 27.90380505526312
$ws.Range("J24").Value = 9.975533488419545
$ws.Range("K24").Value = 11.2895739013687
$ws.Range("L24").Value = 11.09101656567718
$ws.Range("O24").Value = 28.6910613237316
# Row 25
$ws.Range("B25").Value = 13.71899948496729
$ws.Range("C25").Value = 6.913679272198976
$ws.Range("D25").Value = 7.929986127447243
$ws.Range("E25").Value = 12.66571731270737
$ws.Range("F25").Value = 37.31999215323981
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 28.07689411168252
$ws.Range("J25").Value = 10.00905701835531
$ws.Range("K25").Value = 10.9651967182061
$ws.Range("L25").Value = 11.07946813859805
$ws.Range("O25").Value = 28.84554074526537
